# Correcting Relevance Markers Appenzeller-Herzog (2019) - van Dis (2020)
# Update computed metrics for row 3 (file_name = metrics_sim_with_priors.json)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("H3").Value = 0.862173546756529
$ws.Range("I3").Value = 0.02551172074061491
$ws.Range("K3").Value = 174.6875

$ws.Range("Q3").Value = 17
$ws.Range("R3").Value = 25
$ws.Range("S3").Value = 101
$ws.Range("T3").Value = 182
$ws.Range("U3").Value = 231

$ws.Range("V3").Value = 5870
$ws.Range("W3").Value = 5862
$ws.Range("X3").Value = 5786
$ws.Range("Y3").Value = 5705
$ws.Range("Z3").Value = 5656

$ws.Range("AF3").Value = 0.997112
$ws.Range("AG3").Value = 0.995753
$ws.Range("AH3").Value = 0.9828440000000001
$ws.Range("AI3").Value = 0.9690839999999999
$ws.Range("AJ3").Value = 0.960761
